$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# --- "Rectangle 4" (shape id 5): title text "CONCLUSIONS &" -> "CONCLUSION &" ---
$shTitle = $s.Shapes.Item("Rectangle 4")
$shTitle.TextFrame.TextRange.Runs(1).Text = "CONCLUSION &"
# Editing the run text triggers this spAutoFit textbox's height to
# relayout; the diff shows no size change for this shape, so restore the
# original height explicitly.
$shTitle.Height = 128.60882568359375

# --- "Rectangle 14" (shape id 15): bullet text + manual resize ---
$shBullet = $s.Shapes.Item("Rectangle 14")
$shBullet.TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "Low SAT participation rate in both 2017 (2%) & 2018 (3%) "

# Widen/reposition the shape; with spAutoFit this also shrinks its height
# to fit the (now two-line-shorter) text. Values are nudged within the
# COM Single (32-bit float) precision so the EMUs written back land on
# the exact target values.
$shBullet.Left = 375.08195550393697
$shBullet.Top = 168.39740757480317
$shBullet.Width = 563.1265564330708
$shBullet.Height = 68.88621947244094
